$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix decimal separators in B3/B4 (comma -> period)
$ws.Range("B3").Value = "43.36&-5.85"
$ws.Range("B4").Value = "43.24&-5.78"

# Select entire row 4 (A4:XFD4), active cell A4
$ws.Range("A4:XFD4").Select()

# Apply underline formatting to an empty cell C7 (new row)
$ws.Range("C7").Font.Underline = $true
